$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 1857.2084
$ws_ALC.Range("J33").Value = 7219.2
$ws_ALC.Range("L33").Value = 7219.2
$ws_ALC.Range("N33").Value = -7677.2
$ws_ALC.Range("H43").Value = 2733.2222
$ws_ALC.Range("I43").Value = 1574.5
$ws_ALC.Range("J43").Value = 3064.2856
$ws_ALC.Range("K43").Value = 1574.5
$ws_ALC.Range("L43").Value = 3064.2856
$ws_ALC.Range("M43").Value = -1505.5
$ws_ALC.Range("N43").Value = -3202.2856
$ws_ALC.Range("H74").Value = 13032.523
$ws_ALC.Range("I74").Value = 13452.059
$ws_ALC.Range("J74").Value = 11249.5
$ws_ALC.Range("K74").Value = 13452.059
$ws_ALC.Range("L74").Value = 11249.5
$ws_ALC.Range("M74").Value = -12516.059
$ws_ALC.Range("N74").Value = -13121.5
$ws_ALC.Range("H77").Value = 13032.523
$ws_ALC.Range("I77").Value = 13452.059
$ws_ALC.Range("J77").Value = 11249.5
$ws_ALC.Range("K77").Value = 67260.295
$ws_ALC.Range("L77").Value = 56247.5
$ws_ALC.Range("M77").Value = -62580.295
$ws_ALC.Range("N77").Value = -65607.5
$ws_ALC.Range("H115").Value = 498.25
$ws_ALC.Range("J115").Value = 1000
$ws_ALC.Range("L115").Value = 3000
$ws_ALC.Range("N115").Value = -6134
$ws_ALC.Range("H129").Value = 1368.25
$ws_ALC.Range("I129").Value = 828.5
$ws_ALC.Range("J129").Value = 2987.5
$ws_ALC.Range("K129").Value = 2485.5
$ws_ALC.Range("L129").Value = 8962.5
$ws_ALC.Range("M129").Value = 2514.5
$ws_ALC.Range("N129").Value = -18962.5
$ws_ALC.Range("H132").Value = 2082.7817
$ws_ALC.Range("I132").Value = 1595.3396
$ws_ALC.Range("K132").Value = 4786.0188
$ws_ALC.Range("M132").Value = -2256.0188
$ws_ALC.Range("H138").Value = 3128.9814
$ws_ALC.Range("J138").Value = 3413.2703
$ws_ALC.Range("L138").Value = 10239.8109
$ws_ALC.Range("N138").Value = -20519.8109

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 565.5484
$ws_ARM.Range("I2").Value = 565.375
$ws_ARM.Range("K2").Value = 565.375
$ws_ARM.Range("M2").Value = -452.375
$ws_ARM.Range("H32").Value = 9257.469999999999
$ws_ARM.Range("I32").Value = 5426.6665
$ws_ARM.Range("K32").Value = 5426.6665
$ws_ARM.Range("M32").Value = -5139.6665
$ws_ARM.Range("H61").Value = 3647.5908
$ws_ARM.Range("I61").Value = 3421.606
$ws_ARM.Range("J61").Value = 4325.5454
$ws_ARM.Range("K61").Value = 3421.606
$ws_ARM.Range("L61").Value = 4325.5454
$ws_ARM.Range("M61").Value = -3209.606
$ws_ARM.Range("N61").Value = -4749.5454
$ws_ARM.Range("H116").Value = 565.5484
$ws_ARM.Range("I116").Value = 565.375
$ws_ARM.Range("K116").Value = 565.375
$ws_ARM.Range("M116").Value = 1728.625
$ws_ARM.Range("H136").Value = 3647.5908
$ws_ARM.Range("I136").Value = 3421.606
$ws_ARM.Range("J136").Value = 4325.5454
$ws_ARM.Range("K136").Value = 10264.818
$ws_ARM.Range("L136").Value = 12976.6362
$ws_ARM.Range("M136").Value = -7714.818000000001
$ws_ARM.Range("N136").Value = -18076.6362

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 565.5484
$ws_BSM.Range("I3").Value = 565.375
$ws_BSM.Range("K3").Value = 565.375
$ws_BSM.Range("M3").Value = -451.375
$ws_BSM.Range("H80").Value = 22383.611
$ws_BSM.Range("I80").Value = 100039
$ws_BSM.Range("J80").Value = 6852.533
$ws_BSM.Range("K80").Value = 100039
$ws_BSM.Range("L80").Value = 6852.533
$ws_BSM.Range("M80").Value = -99041
$ws_BSM.Range("N80").Value = -8848.532999999999
$ws_BSM.Range("H83").Value = 22383.611
$ws_BSM.Range("I83").Value = 100039
$ws_BSM.Range("J83").Value = 6852.533
$ws_BSM.Range("K83").Value = 500195
$ws_BSM.Range("L83").Value = 34262.665
$ws_BSM.Range("M83").Value = -495203
$ws_BSM.Range("N83").Value = -44246.665
$ws_BSM.Range("H99").Value = 1432.1428
$ws_BSM.Range("I99").Value = 1705.2222
$ws_BSM.Range("J99").Value = 940.6
$ws_BSM.Range("K99").Value = 1705.2222
$ws_BSM.Range("L99").Value = 940.6
$ws_BSM.Range("M99").Value = -207.2221999999999
$ws_BSM.Range("N99").Value = -3936.6
$ws_BSM.Range("H105").Value = 2865086
$ws_BSM.Range("I105").Value = 3340433.8
$ws_BSM.Range("K105").Value = 3340433.8
$ws_BSM.Range("M105").Value = -3338686.8
$ws_BSM.Range("H107").Value = 1945.2354
$ws_BSM.Range("I107").Value = 1324.6364
$ws_BSM.Range("K107").Value = 1324.6364
$ws_BSM.Range("M107").Value = 595.3635999999999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 5155.2856
$ws_CRP.Range("I22").Value = 4348.1665
$ws_CRP.Range("J22").Value = 9998
$ws_CRP.Range("K22").Value = 4348.1665
$ws_CRP.Range("L22").Value = 9998
$ws_CRP.Range("M22").Value = -3998.1665
$ws_CRP.Range("N22").Value = -10698
$ws_CRP.Range("H31").Value = 8606.6875
$ws_CRP.Range("I31").Value = 5876.5806
$ws_CRP.Range("K31").Value = 5876.5806
$ws_CRP.Range("M31").Value = -5581.5806
$ws_CRP.Range("H34").Value = 8606.6875
$ws_CRP.Range("I34").Value = 5876.5806
$ws_CRP.Range("K34").Value = 5876.5806
$ws_CRP.Range("M34").Value = -5674.5806
$ws_CRP.Range("I86").Value = 5004623
$ws_CRP.Range("J86").Value = 5464.7144
$ws_CRP.Range("K86").Value = 5004623
$ws_CRP.Range("L86").Value = 5464.7144
$ws_CRP.Range("M86").Value = -5003500
$ws_CRP.Range("N86").Value = -7710.7144
$ws_CRP.Range("I89").Value = 5004623
$ws_CRP.Range("J89").Value = 5464.7144
$ws_CRP.Range("K89").Value = 25023115
$ws_CRP.Range("L89").Value = 27323.572
$ws_CRP.Range("M89").Value = -25017499
$ws_CRP.Range("N89").Value = -38555.572
$ws_CRP.Range("H94").Value = 1835.1818
$ws_CRP.Range("I94").Value = 771.3333
$ws_CRP.Range("J94").Value = 2234.125
$ws_CRP.Range("K94").Value = 771.3333
$ws_CRP.Range("L94").Value = 2234.125
$ws_CRP.Range("M94").Value = -320.3333
$ws_CRP.Range("N94").Value = -3136.125
$ws_CRP.Range("H105").Value = 3163.4
$ws_CRP.Range("I105").Value = 2626
$ws_CRP.Range("K105").Value = 2626
$ws_CRP.Range("M105").Value = -879
$ws_CRP.Range("H134").Value = 4207.95
$ws_CRP.Range("I134").Value = 3482.0527
$ws_CRP.Range("K134").Value = 10446.1581
$ws_CRP.Range("M134").Value = -7911.158100000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H123").Value = 875
$ws_CUL.Range("I123").Value = 875
$ws_CUL.Range("K123").Value = 2625
$ws_CUL.Range("M123").Value = -175
$ws_CUL.Range("H129").Value = 3692.25
$ws_CUL.Range("J129").Value = 7505.1665
$ws_CUL.Range("L129").Value = 22515.4995
$ws_CUL.Range("N129").Value = -32515.4995

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 402.90625
$ws_GSM.Range("I97").Value = 367.29166
$ws_GSM.Range("J97").Value = 509.75
$ws_GSM.Range("K97").Value = 367.29166
$ws_GSM.Range("L97").Value = 509.75
$ws_GSM.Range("M97").Value = 128.70834
$ws_GSM.Range("N97").Value = -1501.75
$ws_GSM.Range("H105").Value = 62799.332
$ws_GSM.Range("J105").Value = 62799.332
$ws_GSM.Range("L105").Value = 62799.332
$ws_GSM.Range("N105").Value = -69787.33199999999
$ws_GSM.Range("H113").Value = 10317
$ws_GSM.Range("I113").Value = 2917.7273
$ws_GSM.Range("K113").Value = 2917.7273
$ws_GSM.Range("M113").Value = -747.7273
$ws_GSM.Range("H132").Value = 7060.9443
$ws_GSM.Range("I132").Value = 4765.6
$ws_GSM.Range("K132").Value = 14296.8
$ws_GSM.Range("M132").Value = -11766.8

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H36").Value = 24674.092
$ws_LTW.Range("I36").Value = 19962.5
$ws_LTW.Range("J36").Value = 37238.332
$ws_LTW.Range("K36").Value = 19962.5
$ws_LTW.Range("L36").Value = 37238.332
$ws_LTW.Range("M36").Value = -19400.5
$ws_LTW.Range("N36").Value = -38362.332
$ws_LTW.Range("H40").Value = 7593.4116
$ws_LTW.Range("I40").Value = 5517.3335
$ws_LTW.Range("J40").Value = 8725.817999999999
$ws_LTW.Range("K40").Value = 5517.3335
$ws_LTW.Range("L40").Value = 8725.817999999999
$ws_LTW.Range("M40").Value = -5381.3335
$ws_LTW.Range("N40").Value = -8997.817999999999
$ws_LTW.Range("H132").Value = 6210.778
$ws_LTW.Range("I132").Value = 4737.125
$ws_LTW.Range("K132").Value = 14211.375
$ws_LTW.Range("M132").Value = -11681.375

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H122").Value = 3466.7646
$ws_WVR.Range("I122").Value = 3681.5
$ws_WVR.Range("J122").Value = 2464.6667
$ws_WVR.Range("K122").Value = 11044.5
$ws_WVR.Range("L122").Value = 7394.000100000001
$ws_WVR.Range("M122").Value = -8594.5
$ws_WVR.Range("N122").Value = -12294.0001
$ws_WVR.Range("H132").Value = 4113.706
$ws_WVR.Range("I132").Value = 3261.4375
$ws_WVR.Range("J132").Value = 17750
$ws_WVR.Range("K132").Value = 9784.3125
$ws_WVR.Range("L132").Value = 53250
$ws_WVR.Range("M132").Value = -7254.3125
$ws_WVR.Range("N132").Value = -58310
$ws_WVR.Range("H133").Value = 0
$ws_WVR.Range("J133").Value = 0
$ws_WVR.Range("L133").Value = 0

$ws_WVR.Range("N133").ClearContents()

Write-Host "Applied all Moogle_Profits updates"